$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.050941467285156
$ws.Range("B1").Value = 2.520928144454956
$ws.Range("C1").Value = 2.623730897903442
$ws.Range("D1").Value = 3.281381845474243
$ws.Range("E1").Value = 0.9419558644294739
